$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# 9x39 AP and DMG boost: bump base game dmg (G) and IRL Joules (H)
# for ammo_9x39_pab9 (row 21) and ammo_9x39_ap (row 22).
$ws.Range("G21").Value = 0.34
$ws.Range("H21").Value = 1.04

$ws.Range("G22").Value = 0.55
$ws.Range("H22").Value = 1.04

# Update the active cell selection to match the saved view state.
$ws.Activate()
$ws.Range("H20").Select()
